# Updates cryptos list price/volume (and one coin-order swap) to match the
# latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain-looking decimals (e.g. "241.84") as
# text, same as the multi-dot ones (e.g. "29.308.02"). Force NumberFormat
# to Text first so Excel doesn't silently convert these into numeric
# values when we set them below.
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row-by-row value updates (Price in D, Volume(1h) in E; rows 35/36 also
# swap which coin -- ARBITRUM / ImmutableX -- occupies that rank).
$ws.Range("D2").Value = '29.308.02'
$ws.Range("E2").Value = '  -0.24%  '

$ws.Range("D3").Value = '1.872.72'
$ws.Range("E3").Value = '  -0.28%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("E5").Value = '  -0.77%  '

$ws.Range("D6").Value = '241.84'

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").Value = '0.07791'
$ws.Range("E8").Value = '  +0.94%  '

$ws.Range("D9").Value = '0.3106'
$ws.Range("E9").Value = '  -0.37%  '

$ws.Range("D10").Value = '25.07'
$ws.Range("E10").Value = '  -0.50%  '

$ws.Range("D11").Value = '0.08378'
$ws.Range("E11").Value = '  +0.00%  '

$ws.Range("D12").Value = '1.859.62'
$ws.Range("E12").Value = '  -2.75%  '

$ws.Range("D13").Value = '5.236'
$ws.Range("E13").Value = '  -0.12%  '

$ws.Range("D14").Value = '0.7165'
$ws.Range("E14").Value = '  +0.16%  '

$ws.Range("D15").Value = '91.32'
$ws.Range("E15").Value = '  -0.54%  '

$ws.Range("D16").Value = '0.000008396'
$ws.Range("E16").Value = '  +1.25%  '

$ws.Range("D17").Value = '6.140'
$ws.Range("E17").Value = '  +2.75%  '

$ws.Range("D18").Value = '29.324.81'
$ws.Range("E18").Value = '  -0.17%  '

$ws.Range("D19").Value = '240.45'
$ws.Range("E19").Value = '  -1.28%  '

$ws.Range("D20").Value = '2.135.71'
$ws.Range("E20").Value = '  +0.00%  '

$ws.Range("E21").Value = '  -0.35%  '

$ws.Range("D23").Value = '7.741'
$ws.Range("E23").Value = '  -1.85%  '

$ws.Range("E24").Value = '  +0.09%  '

$ws.Range("D25").Value = '0.1597'
$ws.Range("E25").Value = '  -1.44%  '

$ws.Range("D26").Value = '162.46'
$ws.Range("E26").Value = '  -0.83%  '

$ws.Range("D27").Value = '9.032'
$ws.Range("E27").Value = '  -0.04%  '

$ws.Range("D28").Value = '18.49'
$ws.Range("E28").Value = '  -0.41%  '

$ws.Range("E29").Value = '  -0.05%  '

$ws.Range("D30").Value = '4.411'
$ws.Range("E30").Value = '  +0.05%  '

$ws.Range("D31").Value = '4.338'
$ws.Range("E31").Value = '  +0.26%  '

$ws.Range("D32").Value = '1.225'
$ws.Range("E32").Value = '  -5.04%  '

$ws.Range("D33").Value = '0.05356'
$ws.Range("E33").Value = '  +1.95%  '

$ws.Range("D34").Value = '1.940'
$ws.Range("E34").Value = '  +0.61%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.173'
$ws.Range("E35").Value = '  -0.41%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.7461'
$ws.Range("E36").Value = '  -1.39%  '

$ws.Range("E37").Value = '  +0.22%  '

$ws.Range("D38").Value = '0.01874'
$ws.Range("E38").Value = '  +0.54%  '

$ws.Range("D39").Value = '1.241.04'
$ws.Range("E39").Value = '  +6.60%  '

$ws.Range("E40").Value = '  +0.38%  '

$ws.Range("D41").Value = '6.519'
$ws.Range("E41").Value = '  +2.47%  '

$ws.Range("D42").Value = '0.8910'
$ws.Range("E42").Value = '  +0.17%  '

$ws.Range("D43").Value = '109.75'
$ws.Range("E43").Value = '  +4.75%  '

$ws.Range("D44").Value = '72.18'
$ws.Range("E44").Value = '  -2.12%  '

$ws.Range("D46").Value = '2.019.78'
$ws.Range("E46").Value = '  -0.60%  '

$ws.Range("D47").Value = '0.5199'
$ws.Range("E47").Value = '  +0.01%  '

$ws.Range("D48").Value = '1.792'
$ws.Range("E48").Value = '  -0.43%  '

$ws.Range("D49").Value = '9.437'
$ws.Range("E49").Value = '  +0.29%  '

$ws.Range("D50").Value = '0.4333'
$ws.Range("E50").Value = '  +0.50%  '

$ws.Range("D51").Value = '7.089'
$ws.Range("E51").Value = '  +0.52%  '
